$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction: header cell B2 should read "total" (was mistakenly
# "unnamed: 1_level_1" from the raw pandas export).
$ws.Range("B2").Value = "total"

# Remove the two empty sub-heading rows that only carried a stray
# label with no data ("situação do domicílio" and "grandes regiões"),
# shifting the rows below them up.
$ws.Rows.Item(8).EntireRow.Delete()
$ws.Rows.Item(5).EntireRow.Delete()
